$d = $word.ActiveDocument

# Add a new paragraph at the end of the document, then insert a hyperlink
# pointing at "This is a test.pptx" with visible text "This is a test.pptx"
$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()
$end.Collapse(0)

$d.Hyperlinks.Add($end, "This is a test.pptx", $null, $null, "This is a test.pptx")
